# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to keep numeric-looking strings (e.g. "9.00",
# "0.0000171", "60.957.72") as text instead of auto-converting them to numbers,
# matching how the Price column values are stored in the source data.

$ws.Range('D2').Value = '''60.957.72'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '''3.395.44'
$ws.Range('E3').Value = '  -1.09%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''570.37'
$ws.Range('D6').Value = '''142.01'
$ws.Range('E6').Value = '  -2.37%  '
$ws.Range('D7').Value = '''3.395.67'
$ws.Range('E7').Value = '  -1.10%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.475'
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('D10').Value = '''7.54'
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('D12').Value = '''0.395'
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('D13').Value = '''3.974.50'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('D14').Value = '''28.32'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '''0.0000171'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '''3.396.87'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').Value = '''61.027.50'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = '''6.19'
$ws.Range('E19').Value = '  -1.20%  '
$ws.Range('D20').Value = '''13.91'
$ws.Range('E20').Value = '  -1.89%  '
$ws.Range('D21').Value = '''9.00'
$ws.Range('E21').Value = '  -4.19%  '
$ws.Range('D22').Value = '''386.05'
$ws.Range('E22').Value = '  -2.39%  '
$ws.Range('D23').Value = '''0.558'
$ws.Range('E23').Value = '  -1.29%  '
$ws.Range('D24').Value = '''74.10'
$ws.Range('E24').Value = '  +1.52%  '
$ws.Range('E25').Value = '  +0.58%  '
$ws.Range('E26').Value = '  -4.31%  '
$ws.Range('D27').Value = '''3.534.32'
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('E28').Value = '  +0.24%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  -2.65%  '
$ws.Range('E31').Value = '  -2.67%  '
$ws.Range('E32').Value = '  -1.73%  '
$ws.Range('E33').Value = '  -2.75%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '''23.57'
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('D36').Value = '''6.99'
$ws.Range('E36').Value = '  -0.40%  '
$ws.Range('D37').Value = '''167.32'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').Value = '''3.425.60'
$ws.Range('E38').Value = '  -1.00%  '
$ws.Range('E39').Value = '  -2.49%  '
$ws.Range('E40').Value = '  -4.66%  '
$ws.Range('D41').Value = '''28.16'
$ws.Range('E41').Value = '  +4.76%  '
$ws.Range('D42').Value = '''0.0775'
$ws.Range('E42').Value = '  -1.45%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('E44').Value = '  -2.31%  '
$ws.Range('D45').Value = '''42.12'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').Value = '  -1.36%  '
$ws.Range('D47').Value = '''1.68'
$ws.Range('E47').Value = '  -3.35%  '
$ws.Range('D48').Value = '''1.14'
$ws.Range('E48').Value = '  -2.29%  '
$ws.Range('D49').Value = '''2.486.42'
$ws.Range('E49').Value = '  -3.61%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '''6.83'
$ws.Range('E50').Value = '  -1.34%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '''23.11'
$ws.Range('E51').Value = '  -0.55%  '
